$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J: header "URL_AvenueDeLaBrique", formatted like the other
# header cells (bold / centered / bordered -- same style as A1:I1).
$ws.Range("J1").Value = "URL_AvenueDeLaBrique"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row for set 77243 (Voiture F1(R) Oracle Red Bull Racing RB20).
$r = 18
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = "77243"
$ws.Cells.Item($r, 2).Value = "Voiture F1® Oracle Red Bull Racing RB20"
$ws.Cells.Item($r, 3).Value = "N/A"
$ws.Cells.Item($r, 4).Value = "Speed Champions"
$ws.Cells.Item($r, 5).Value = "https://www.lego.com/cdn/cs/set/assets/blt8dac22afe99a2c70/77243_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Cells.Item($r, 6).Value = "https://www.lego.com/fr-fr/product/77243"
$ws.Cells.Item($r, 10).Value = "https://www.avenuedelabrique.com/lego-speed-champions/77243-voiture-f1-oracle-red-bull-racing-rb20/p10441"
